# Generate Report for Handback
#
# - Status text changes from "Ready for handoff" to "Handed back: in sync
#   with en-US" everywhere it appears (Overview + both language sheets).
# - Each language sheet (zh-cn, de-de) gets two new hyperlinked columns
#   filled in for every data row: F = "Latest Target File" (the source
#   .md, same link as column A) and G = "Latest Handback File" (the
#   handed-back .xlf, same link as column D).
# - Column H ("Latest Handback DateTime") is stamped with the real
#   handback time instead of the zero-date placeholder.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet: just the status text -------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($addr in @("B2", "C2", "B3", "C3")) {
    $cell = $wsOverview.Range($addr)
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

# ---- Per-language sheets --------------------------------------------------
# Row -> (source .md display/target, handback .xlf display/target, handback datetime)
$langSheets = @("zh-cn", "de-de")

$rowInfo = @{
    2 = @{
        MdDisplay = "50527144-a073-47db-9c7c-0e38a0676b0f.md"
    }
    3 = @{
        MdDisplay = "a45f4c37-42f6-490b-8d01-a84c223ce2ca.md"
    }
}

$handbackDatetime = @{
    "zh-cn" = "2016-03-24 10:16:09"
    "de-de" = "2016-03-24 10:16:16"
}

foreach ($sheetName in $langSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in @(2, 3)) {
        # Status column C
        $statusCell = $ws.Range("C$row")
        if ($statusCell.Value2 -eq $oldStatus) {
            $statusCell.Value = $newStatus
        }

        # Existing hyperlinked cells we mirror into F/G
        $aCell = $ws.Range("A$row")
        $dCell = $ws.Range("D$row")

        $mdDisplay = $aCell.Text
        $xlfDisplay = $dCell.Text

        $aLink = $ws.Hyperlinks.Item(1)
        # Find the actual hyperlink objects bound to A$row / D$row so we can
        # reuse their target addresses for the new F/G hyperlinks.
        $aAddress = $null
        $dAddress = $null
        foreach ($hl in $ws.Hyperlinks) {
            if ($hl.Range.Address -eq $aCell.Address) { $aAddress = $hl.Address }
            if ($hl.Range.Address -eq $dCell.Address) { $dAddress = $hl.Address }
        }

        # F column: "Latest Target File" - same file/link as the source .md
        $fCell = $ws.Range("F$row")
        $ws.Hyperlinks.Add($fCell, $aAddress, "", "", $mdDisplay)
        $fCell.Font.Underline = $true
        $fCell.Font.Color = $aCell.Font.Color

        # G column: "Latest Handback File" - same file/link as the handoff .xlf
        $gCell = $ws.Range("G$row")
        $ws.Hyperlinks.Add($gCell, $dAddress, "", "", $xlfDisplay)
        $gCell.Font.Underline = $true
        $gCell.Font.Color = $dCell.Font.Color

        # H column: "Latest Handback DateTime"
        $hCell = $ws.Range("H$row")
        $hCell.Value = $handbackDatetime[$sheetName]
    }
}
